$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current "season_winter" row (row 12) to make
# room for the new hotspot-analyze / investigate / season-title strings.
$ws.Rows.Item(12).Resize(5).Insert()

# New Key/Value pairs introduced by this change (rows 12-16).
$newEntries = @(
    @("hotspotAnalyze_title", "Atmospheric Reading"),
    @("analyzing", "ANALYZING"),
    @("incompatible", "INCOMPATIBLE!"),
    @("investigate", "INVESTIGATE"),
    @("season_title", "Season")
)

$row = 12
foreach ($entry in $newEntries) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

# Update the remembered selection to match the authored workbook.
$ws.Range("B14").Select()
